$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The subject_label values in column A for rows 2 and 3 change from the
# text "012FX17968" to the numeric value 175055.
$ws.Range("A2").Value = 175055
$ws.Range("A3").Value = 175055
